$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.132.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2831"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06458"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("E11").Value = "  -3.21%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07729"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.71%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.79%  "

# Row 14
$ws.Range("E14").Value = "  -3.58%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6800"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.057"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "265.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.096.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007554"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.106.12"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.156"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.098"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.288"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.885"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.366"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09834"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.450"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.204"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.976"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.20%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04654"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6841"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "

# Row 38
$ws.Range("E38").Value = "  +0.35%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01809"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.715"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.34%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.278"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.0000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.880"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.95%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4045"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.90%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.061"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.03%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "924.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.932"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.08%  "
